$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: PTP_checkbox (J4) unchecked, amended_k1 (L4) checked
$ws.Range("J4").Value = $false
$ws.Range("L4").Value = $true

# Update active cell selection to J5
$ws.Range("J5").Select()
